$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting existing rows 103:147 down to 104:148
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new data point
$ws.Cells.Item(103, 1).Value = 1
$ws.Cells.Item(103, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(103, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(103, 4).Value = 45229
$ws.Cells.Item(103, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(103, 5).Value = 15
$ws.Cells.Item(103, 6).Value = 100112038
$ws.Cells.Item(103, 7).Value = "Cebollín baby"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 170
$ws.Cells.Item(103, 11).Value = 1800
$ws.Cells.Item(103, 12).Value = 2000
$ws.Cells.Item(103, 13).Value = 1906
$ws.Cells.Item(103, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(103, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(103, 16).Value = 953
$ws.Cells.Item(103, 17).Value = 2
$ws.Cells.Item(103, 18).Value = "Hortaliza"
